$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1669.2623
$ws.Range("J17").Value = 1669.2623
$ws.Range("L17").Value = 5007.7869
$ws.Range("N17").Value = -5343.7869

$ws.Range("H61").Value = 649.1429000000001
$ws.Range("I61").Value = 665.6667
$ws.Range("J61").Value = 550
$ws.Range("K61").Value = 1997.0001
$ws.Range("L61").Value = 1650
$ws.Range("M61").Value = -1825.0001
$ws.Range("N61").Value = -1994

$ws.Range("H64").Value = 4999.5
$ws.Range("I64").Value = 4999.6665
$ws.Range("K64").Value = 4999.6665
$ws.Range("M64").Value = -4751.6665

$ws.Range("H67").Value = 4999.5
$ws.Range("I67").Value = 4999.6665
$ws.Range("K67").Value = 4999.6665
$ws.Range("M67").Value = -4141.6665

$ws.Range("H70").Value = 2714.0715
$ws.Range("I70").Value = 6299.75
$ws.Range("J70").Value = 1279.8
$ws.Range("K70").Value = 18899.25
$ws.Range("L70").Value = 3839.4
$ws.Range("M70").Value = -18629.25
$ws.Range("N70").Value = -4379.4

$ws.Range("H73").Value = 2714.0715
$ws.Range("I73").Value = 6299.75
$ws.Range("J73").Value = 1279.8
$ws.Range("K73").Value = 18899.25
$ws.Range("L73").Value = 3839.4
$ws.Range("M73").Value = -17963.25
$ws.Range("N73").Value = -5711.4

$ws.Range("H74").Value = 3692.75
$ws.Range("I74").Value = 3692.75
$ws.Range("K74").Value = 3692.75
$ws.Range("M74").Value = -2756.75

$ws.Range("H75").Value = 48571.285
$ws.Range("I75").Value = 25000
$ws.Range("J75").Value = 52499.832
$ws.Range("K75").Value = 25000
$ws.Range("L75").Value = 52499.832
$ws.Range("N75").Value = -54371.832
$ws.Range("M75").Value = -24064

$ws.Range("H77").Value = 3692.75
$ws.Range("I77").Value = 3692.75
$ws.Range("K77").Value = 18463.75
$ws.Range("M77").Value = -13783.75

$ws.Range("H78").Value = 48571.285
$ws.Range("I78").Value = 25000
$ws.Range("J78").Value = 52499.832
$ws.Range("K78").Value = 75000
$ws.Range("L78").Value = 157499.496
$ws.Range("N78").Value = -166859.496
$ws.Range("M78").Value = -70320

$ws.Range("H101").Value = 289
$ws.Range("I101").Value = 289
$ws.Range("K101").Value = 867
$ws.Range("M101").Value = 755

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H131").Value = 1650
$ws.Range("I131").Value = 1650
$ws.Range("K131").Value = 4950
$ws.Range("M131").Value = 90

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1971.1428
$ws.Range("I2").Value = 1799.6
$ws.Range("K2").Value = 1799.6
$ws.Range("M2").Value = -1686.6

$ws.Range("H4").Value = 930.5
$ws.Range("I4").Value = 930.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 930.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -814.5
$ws.Range("N4").ClearContents()

$ws.Range("H32").Value = 2978.8928
$ws.Range("I32").Value = 2718.889
$ws.Range("K32").Value = 2718.889
$ws.Range("M32").Value = -2431.889

$ws.Range("H61").Value = 2333.3333
$ws.Range("I61").Value = 2333.3333
$ws.Range("K61").Value = 2333.3333
$ws.Range("M61").Value = -2121.3333

$ws.Range("H102").Value = 1598
$ws.Range("I102").Value = 1598
$ws.Range("K102").Value = 1598
$ws.Range("M102").Value = 24

$ws.Range("H116").Value = 1971.1428
$ws.Range("I116").Value = 1799.6
$ws.Range("K116").Value = 1799.6
$ws.Range("M116").Value = 494.4000000000001

$ws.Range("H132").Value = 1357.091
$ws.Range("I132").Value = 1357.091
$ws.Range("K132").Value = 4071.273
$ws.Range("M132").Value = -1541.273

$ws.Range("H136").Value = 2333.3333
$ws.Range("I136").Value = 2333.3333
$ws.Range("K136").Value = 6999.999899999999
$ws.Range("M136").Value = -4449.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1971.1428
$ws.Range("I3").Value = 1799.6
$ws.Range("K3").Value = 1799.6
$ws.Range("M3").Value = -1685.6

$ws.Range("H86").Value = 3057.6875
$ws.Range("I86").Value = 3057.6875
$ws.Range("K86").Value = 3057.6875
$ws.Range("M86").Value = -1934.6875

$ws.Range("H89").Value = 3057.6875
$ws.Range("I89").Value = 3057.6875
$ws.Range("K89").Value = 15288.4375
$ws.Range("M89").Value = -9672.4375

$ws.Range("H105").Value = 2633.6316
$ws.Range("I105").Value = 2025.7307
$ws.Range("K105").Value = 2025.7307
$ws.Range("M105").Value = -278.7307000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4080
$ws.Range("I58").Value = 4725
$ws.Range("J58").Value = 1500
$ws.Range("K58").Value = 4725
$ws.Range("L58").Value = 1500
$ws.Range("M58").Value = -4522
$ws.Range("N58").Value = -1906

$ws.Range("H86").Value = 8000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 8000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H94").Value = 1188.8
$ws.Range("J94").Value = 698
$ws.Range("L94").Value = 698
$ws.Range("N94").Value = -1600

$ws.Range("H136").Value = 4080
$ws.Range("I136").Value = 4725
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 14175
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -11625
$ws.Range("N136").Value = -9600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1592.8695
$ws.Range("J26").Value = 859.75
$ws.Range("L26").Value = 2579.25
$ws.Range("N26").Value = -3155.25

$ws.Range("H44").Value = 2474.5
$ws.Range("I44").Value = 4000
$ws.Range("J44").Value = 949
$ws.Range("K44").Value = 12000
$ws.Range("L44").Value = 2847
$ws.Range("M44").Value = -11602
$ws.Range("N44").Value = -3643

$ws.Range("H68").Value = 623.3333
$ws.Range("I68").Value = 598
$ws.Range("K68").Value = 1794
$ws.Range("M68").Value = -983

$ws.Range("H71").Value = 623.3333
$ws.Range("I71").Value = 598
$ws.Range("K71").Value = 5382
$ws.Range("M71").Value = -1326

$ws.Range("H131").Value = 2309.2727
$ws.Range("I131").Value = 1857.5714
$ws.Range("J131").Value = 3099.75
$ws.Range("K131").Value = 5572.7142
$ws.Range("L131").Value = 9299.25
$ws.Range("M131").Value = -532.7142000000003
$ws.Range("N131").Value = -19379.25

$ws.Range("H132").Value = 1029.125
$ws.Range("I132").Value = 1184.8
$ws.Range("J132").Value = 769.6667
$ws.Range("K132").Value = 10663.2
$ws.Range("L132").Value = 6927.0003
$ws.Range("M132").Value = -8133.199999999999
$ws.Range("N132").Value = -11987.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3614.7778
$ws.Range("I80").Value = 4118
$ws.Range("J80").Value = 2985.75
$ws.Range("K80").Value = 4118
$ws.Range("L80").Value = 2985.75
$ws.Range("M80").Value = -3120
$ws.Range("N80").Value = -4981.75

$ws.Range("H83").Value = 3614.7778
$ws.Range("I83").Value = 4118
$ws.Range("J83").Value = 2985.75
$ws.Range("K83").Value = 20590
$ws.Range("L83").Value = 14928.75
$ws.Range("M83").Value = -15598
$ws.Range("N83").Value = -24912.75

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H132").Value = 2369.476
$ws.Range("I132").Value = 2115.4546
$ws.Range("K132").Value = 6346.3638
$ws.Range("M132").Value = -3816.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1149.8182
$ws.Range("I55").Value = 299.75
$ws.Range("J55").Value = 1635.5714
$ws.Range("K55").Value = 299.75
$ws.Range("L55").Value = 1635.5714
$ws.Range("M55").Value = -126.75
$ws.Range("N55").Value = -1981.5714

$ws.Range("H82").Value = 1079
$ws.Range("I82").Value = 2142.3333
$ws.Range("J82").Value = 547.3333
$ws.Range("K82").Value = 2142.3333
$ws.Range("L82").Value = 547.3333
$ws.Range("M82").Value = -1781.3333
$ws.Range("N82").Value = -1269.3333

$ws.Range("H85").Value = 1079
$ws.Range("I85").Value = 2142.3333
$ws.Range("J85").Value = 547.3333
$ws.Range("K85").Value = 2142.3333
$ws.Range("L85").Value = 547.3333
$ws.Range("M85").Value = -894.3332999999998
$ws.Range("N85").Value = -3043.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H62").Value = 11636.385
$ws.Range("J62").Value = 12567.3
$ws.Range("L62").Value = 12567.3
$ws.Range("N62").Value = -13815.3

$ws.Range("H65").Value = 11636.385
$ws.Range("J65").Value = 12567.3
$ws.Range("L65").Value = 62836.5
$ws.Range("N65").Value = -69076.5

$ws.Range("H81").Value = 4318.1816
$ws.Range("I81").Value = 1900
$ws.Range("J81").Value = 4855.5557
$ws.Range("K81").Value = 3800
$ws.Range("L81").Value = 9711.1114
$ws.Range("M81").Value = -2739
$ws.Range("N81").Value = -11833.1114

$ws.Range("H84").Value = 4318.1816
$ws.Range("I84").Value = 1900
$ws.Range("J84").Value = 4855.5557
$ws.Range("K84").Value = 19000
$ws.Range("L84").Value = 48555.557
$ws.Range("M84").Value = -13696
$ws.Range("N84").Value = -59163.557
